$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 189.0573523333333
$ws.Range("H2").Value = 567.172057
$ws.Range("I2").Value = 0.1182556374491171
$ws.Range("J2").Value = 0.1182556374491171
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.017765
$ws.Range("N2").Value = 0.053295
$ws.Range("O2").Value = 0.005225147533577419
$ws.Range("P2").Value = 0.005225147533577419
$ws.Range("Q2").Value = 3.358603864201667
$ws.Range("R2").Value = 30.227434777815
$ws.Range("S2").Value = 0.0006179031523488795
$ws.Range("T2").Value = 0.0006179031523488795

$ws.Range("G3").Value = 189.0573523333333
$ws.Range("H3").Value = 567.172057
$ws.Range("I3").Value = 0.1182556374491171
$ws.Range("J3").Value = 0.1182556374491171
$ws.Range("M3").Value = 0.8788360000000001
$ws.Range("N3").Value = 2.636508
$ws.Range("O3").Value = 0.258488474968705
$ws.Range("P3").Value = 0.258488474968705
$ws.Range("Q3").Value = 166.1504072952173
$ws.Range("R3").Value = 1495.353665656956
$ws.Range("S3").Value = 0.03056771938067436
$ws.Range("T3").Value = 0.03056771938067436

$ws.Range("G4").Value = 189.0573523333333
$ws.Range("H4").Value = 567.172057
$ws.Range("I4").Value = 0.1182556374491171
$ws.Range("J4").Value = 0.1182556374491171
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.503303
$ws.Range("N4").Value = 7.509909
$ws.Range("O4").Value = 0.7362863774977175
$ws.Range("P4").Value = 0.7362863774977175
$ws.Range("Q4").Value = 473.2678372680903
$ws.Range("R4").Value = 4259.410535412812
$ws.Range("S4").Value = 0.08707001491609384
$ws.Range("T4").Value = 0.08707001491609384

$ws.Range("G5").Value = 930.1503093333332
$ws.Range("I5").Value = 0.5818103152093762
$ws.Range("J5").Value = 0.5818103152093762
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.017765
$ws.Range("N5").Value = 0.053295
$ws.Range("O5").Value = 0.005225147533577419
$ws.Range("P5").Value = 0.005225147533577419
$ws.Range("Q5").Value = 16.52412024530667
$ws.Range("R5").Value = 148.71708220776
$ws.Range("S5").Value = 0.003040044733526173
$ws.Range("T5").Value = 0.003040044733526173

$ws.Range("G6").Value = 930.1503093333332
$ws.Range("I6").Value = 0.5818103152093762
$ws.Range("J6").Value = 0.5818103152093762
$ws.Range("M6").Value = 0.8788360000000001
$ws.Range("N6").Value = 2.636508
$ws.Range("O6").Value = 0.258488474968705
$ws.Range("P6").Value = 0.258488474968705
$ws.Range("Q6").Value = 817.4495772532692
$ws.Range("R6").Value = 7357.046195279424
$ws.Range("S6").Value = 0.1503912610995332
$ws.Range("T6").Value = 0.1503912610995332

$ws.Range("G7").Value = 930.1503093333332
$ws.Range("I7").Value = 0.5818103152093762
$ws.Range("J7").Value = 0.5818103152093762
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.503303
$ws.Range("N7").Value = 7.509909
$ws.Range("O7").Value = 0.7362863774977175
$ws.Range("P7").Value = 0.7362863774977175
$ws.Range("Q7").Value = 2328.448059805061
$ws.Range("R7").Value = 20956.03253824555
$ws.Range("S7").Value = 0.4283790093763168
$ws.Range("T7").Value = 0.4283790093763168

$ws.Range("G8").Value = 420.6651306666666
$ws.Range("H8").Value = 1261.995392
$ws.Range("I8").Value = 0.2631266256807295
$ws.Range("J8").Value = 0.2631266256807295
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.017765
$ws.Range("N8").Value = 0.053295
$ws.Range("O8").Value = 0.005225147533577419
$ws.Range("P8").Value = 0.005225147533577419
$ws.Range("Q8").Value = 7.473116046293332
$ws.Range("R8").Value = 67.25804441663999
$ws.Range("S8").Value = 0.001374875439194212
$ws.Range("T8").Value = 0.001374875439194212

$ws.Range("G9").Value = 420.6651306666666
$ws.Range("H9").Value = 1261.995392
$ws.Range("I9").Value = 0.2631266256807295
$ws.Range("J9").Value = 0.2631266256807295
$ws.Range("M9").Value = 0.8788360000000001
$ws.Range("N9").Value = 2.636508
$ws.Range("O9").Value = 0.258488474968705
$ws.Range("P9").Value = 0.258488474968705
$ws.Range("Q9").Value = 369.6956607745706
$ws.Range("R9").Value = 3327.260946971136
$ws.Range("S9").Value = 0.06801520019587307
$ws.Range("T9").Value = 0.06801520019587308

$ws.Range("G10").Value = 420.6651306666666
$ws.Range("H10").Value = 1261.995392
$ws.Range("I10").Value = 0.2631266256807295
$ws.Range("J10").Value = 0.2631266256807295
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.503303
$ws.Range("N10").Value = 7.509909
$ws.Range("O10").Value = 0.7362863774977175
$ws.Range("P10").Value = 0.7362863774977175
$ws.Range("Q10").Value = 1053.052283593258
$ws.Range("R10").Value = 9477.470552339326
$ws.Range("S10").Value = 0.1937365500456622
$ws.Range("T10").Value = 0.1937365500456622

$ws.Range("G11").Value = 58.84466766666667
$ws.Range("H11").Value = 176.534003
$ws.Range("I11").Value = 0.03680742166077718
$ws.Range("J11").Value = 0.03680742166077718
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.017765
$ws.Range("N11").Value = 0.053295
$ws.Range("O11").Value = 0.005225147533577419
$ws.Range("P11").Value = 0.005225147533577419
$ws.Range("Q11").Value = 1.045375521098333
$ws.Range("R11").Value = 9.408379689885001
$ws.Range("S11").Value = 0.0001923242085081539
$ws.Range("T11").Value = 0.0001923242085081539

$ws.Range("G12").Value = 58.84466766666667
$ws.Range("H12").Value = 176.534003
$ws.Range("I12").Value = 0.03680742166077718
$ws.Range("J12").Value = 0.03680742166077718
$ws.Range("M12").Value = 0.8788360000000001
$ws.Range("N12").Value = 2.636508
$ws.Range("O12").Value = 0.258488474968705
$ws.Range("P12").Value = 0.258488474968705
$ws.Range("Q12").Value = 51.71481235350267
$ws.Range("R12").Value = 465.4333111815241
$ws.Range("S12").Value = 0.009514294292624374
$ws.Range("T12").Value = 0.009514294292624374

$ws.Range("G13").Value = 58.84466766666667
$ws.Range("H13").Value = 176.534003
$ws.Range("I13").Value = 0.03680742166077718
$ws.Range("J13").Value = 0.03680742166077718
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.503303
$ws.Range("N13").Value = 7.509909
$ws.Range("O13").Value = 0.7362863774977175
$ws.Range("P13").Value = 0.7362863774977175
$ws.Range("Q13").Value = 147.3060331039697
$ws.Range("R13").Value = 1325.754297935727
$ws.Range("S13").Value = 0.02710080315964465
$ws.Range("T13").Value = 0.02710080315964465

